$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was inserted before the existing row 149
# (Fecha 2022-02-07 / serial 44582), pushing every following record
# (old rows 149-244) down by one row.
$ws.Rows.Item(149).Insert()

$ws.Cells.Item(149, 1).Value2 = 7
$ws.Cells.Item(149, 2).Value2 = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(149, 3).Value2 = "Ñuble"
$ws.Cells.Item(149, 4).Value2 = 44582
$ws.Cells.Item(149, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(149, 5).Value2 = 16
$ws.Cells.Item(149, 6).Value2 = 100114013
$ws.Cells.Item(149, 7).Value2 = "Zanahoria"
$ws.Cells.Item(149, 8).Value2 = "Sin especificar"
$ws.Cells.Item(149, 9).Value2 = "Primera"
$ws.Cells.Item(149, 10).Value2 = 100
$ws.Cells.Item(149, 11).Value2 = 7000
$ws.Cells.Item(149, 12).Value2 = 7500
$ws.Cells.Item(149, 13).Value2 = 7250
$ws.Cells.Item(149, 14).Value2 = "$/saco 20 kilos"
$ws.Cells.Item(149, 15).Value2 = "Provincia de Diguillín"
$ws.Cells.Item(149, 16).Value2 = 362
$ws.Cells.Item(149, 17).Value2 = 20
$ws.Cells.Item(149, 18).Value2 = "Hortaliza"
